$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled data
$ws.Range("F2").Value = -1
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 3
$ws.Range("F9").Value = -1
$ws.Range("F11").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = -5
$ws.Range("F14").Value = -6
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = -2
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = 1
$ws.Range("F20").Value = -5
$ws.Range("F21").Value = -2
$ws.Range("F22").Value = 1
